# Update automatico via Actualizar 05-09-2020 05-37-52
#
# Adds 7 new hospital/clinic records (rows 129-135) to the HOSPITALES_HN
# table on the HOSPITALES sheet, backfilling the Pais_cod3 / Pais /
# Nivel_Admin columns (B:D) for all the previously-added rows (97-135)
# that were missing them, resizes the table / filter range accordingly,
# and updates the frozen-pane / selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOSPITALES")

# ---------------------------------------------------------------------
# 1. Backfill Pais_cod3 (B) / Pais (C) / Nivel_Admin (D) for rows 97-135.
#    Every data row in this sheet carries "HND" / "Honduras" / 3 in
#    these columns; rows 97-128 (and the new rows below) were missing
#    them.
# ---------------------------------------------------------------------
for ($r = 97; $r -le 135; $r++) {
    $ws.Range("B$r").Value = "HND"
    $ws.Range("C$r").Value = "Honduras"
    $ws.Range("D$r").Value = 3
}

# ---------------------------------------------------------------------
# 2. Row 129 gains its Admin1nombre (G) and Admin2nombre (K) values -
#    the facility itself (columns S/U/V/W) was already present.
# ---------------------------------------------------------------------
$ws.Range("G129").Value = "Lempira"
$ws.Range("K129").Value = "Santa Cruz"

# ---------------------------------------------------------------------
# 3. Brand-new rows 130-135: additional Lempira / Valle facilities.
# ---------------------------------------------------------------------

# Row 130 - Gualcinse, Lempira
$ws.Range("G130").Value = "Lempira"
$ws.Range("K130").Value = "Gualcinse"
$ws.Range("S130").Value = "hospital"
$ws.Range("U130").Value = "Hospital del Sur Dr. Lempira"
$ws.Range("V130").Value = 14.126211
$ws.Range("W130").Value = -88.543222999999998

# Row 131 - Lempira (centro de salud)
$ws.Range("G131").Value = "Lempira"
$ws.Range("S131").Value = "centro de salud"
$ws.Range("U131").Value = "Centro de Salud"
$ws.Range("V131").Value = 14.114020999999999
$ws.Range("W131").Value = -88.651722000000007

# Row 132 - Nacaome, Valle
$ws.Range("G132").Value = "Valle"
$ws.Range("K132").Value = "Nacaome"
$ws.Range("S132").Value = "clínica"
$ws.Range("U132").Value = "Policlínica Nacaome"
$ws.Range("V132").Value = 13.530787999999999
$ws.Range("W132").Value = -87.498217999999994

# Row 133 - Nacaome, Valle
$ws.Range("G133").Value = "Valle"
$ws.Range("K133").Value = "Nacaome"
$ws.Range("S133").Value = "hospital"
$ws.Range("U133").Value = "Centro Médico Juárez"
$ws.Range("V133").Value = 13.532709000000001
$ws.Range("W133").Value = -87.492543999999995

# Row 134 - Langue, Valle
$ws.Range("G134").Value = "Valle"
$ws.Range("K134").Value = "Langue"
$ws.Range("S134").Value = "centro de salud"
$ws.Range("U134").Value = "Centro de Salud de Langue"
$ws.Range("V134").Value = 13.62003
$ws.Range("W134").Value = -87.657388999999995

# Row 135 - Guascorán, Valle
$ws.Range("G135").Value = "Valle"
$ws.Range("K135").Value = "Guascorán"
$ws.Range("S135").Value = "clínica"
$ws.Range("U135").Value = "Medicenter"
$ws.Range("V135").Value = 13.610545
$ws.Range("W135").Value = -87.752651

# ---------------------------------------------------------------------
# 4. Grow the HOSPITALES_HN table / AutoFilter from A1:W131 to A1:W139
#    and keep the _FilterDatabase defined name in sync.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("HOSPITALES_HN")
$lo.Resize($ws.Range("A1:W139"))

$filterName = $wb.Names.Item("HOSPITALES!_FilterDatabase")
$filterName.RefersTo = "=HOSPITALES!`$A`$1:`$W`$139"

# ---------------------------------------------------------------------
# 5. Refresh the frozen header pane / current selection to match the
#    newly-scrolled view over the new rows.
# ---------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D96:D135").Select()
